$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aCell = $ws.Range("A312")
$aCell.NumberFormat = "@"
$aCell.Value = "000311"
$aCell.Style = "Normal"
$ws.Range("B312").Value = "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-May-2023 16:46:41"

$aCell = $ws.Range("A313")
$aCell.NumberFormat = "@"
$aCell.Value = "000312"
$aCell.Style = "Normal"
$ws.Range("B313").Value = "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-May-2023 17:02:54"

$aCell = $ws.Range("A314")
$aCell.NumberFormat = "@"
$aCell.Value = "000313"
$aCell.Style = "Normal"
$ws.Range("B314").Value = "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-May-2023 17:04:14"

$aCell = $ws.Range("A315")
$aCell.NumberFormat = "@"
$aCell.Value = "000314"
$aCell.Style = "Normal"
$ws.Range("B315").Value = "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Read_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-May-2023 17:05:24"
